# Update the "want to go" counts (column F) on the "展览" and "全部类型"
# sheets to reflect newly generated output numbers.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value  = 44
$wsExhibit.Range("F4").Value  = 1548
$wsExhibit.Range("F5").Value  = 241
$wsExhibit.Range("F6").Value  = 53
$wsExhibit.Range("F7").Value  = 998
$wsExhibit.Range("F8").Value  = 10081
$wsExhibit.Range("F11").Value = 250
$wsExhibit.Range("F14").Value = 6987
$wsExhibit.Range("F16").Value = 654
$wsExhibit.Range("F18").Value = 217

# Sheet "全部类型"
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 44
$wsAll.Range("F4").Value  = 1548
$wsAll.Range("F5").Value  = 241
$wsAll.Range("F7").Value  = 53
$wsAll.Range("F8").Value  = 998
$wsAll.Range("F11").Value = 10081
$wsAll.Range("F14").Value = 250
$wsAll.Range("F17").Value = 6987
$wsAll.Range("F19").Value = 654
$wsAll.Range("F21").Value = 217
